$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$sh.TextFrame.TextRange.Text = "11/5/2019"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "CSE 4361 Lecture 17") {
                $tr.Text = "CSE 4361 Lecture 18"
            }
        }
    }
}
Write-Output "done"
